$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new column at E. This shifts the helper lookup data that used to
# live in columns F/G one column to the right (F->G, G->H), exactly as a
# real "insert column" does in Excel.
# ---------------------------------------------------------------------------
$ws.Columns("E:E").Insert()
$ws.Columns("E:E").ColumnWidth = 13.140625

# Copy the formatting (borders/number format/font) of the "Độ ưu tiên" data
# column (D) onto the brand-new column (E) so the new column visually
# matches the rest of the table.
$ws.Range("D3:D48").Copy()
$ws.Range("E3:E48").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Grow Table1 so it covers the new column, then name it.
# ---------------------------------------------------------------------------
$t = $ws.ListObjects.Item(1)
$t.Resize($ws.Range("A3:E48"))
$ws.Range("E3").Value = "Hoàn thành"

# ---------------------------------------------------------------------------
# Fill in the "Hoàn thành" (completion) percentages for the rows that have
# one; the rest of the column is left blank.
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = 0.9
$ws.Range("E5").Value = 0.9
$ws.Range("E9").Value = 0.9
$ws.Range("E10").Value = 0.9
$ws.Range("E21").Value = 1
$ws.Range("E22").Value = 0.5

# Whole data column is formatted as a percentage.
$ws.Range("E4:E48").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# Re-merge the title cell across the (now wider) table width.
# ---------------------------------------------------------------------------
$ws.Range("A1:D2").UnMerge()
$ws.Range("A1:E2").Merge()

# ---------------------------------------------------------------------------
# Match the selection left behind by the editor.
# ---------------------------------------------------------------------------
$ws.Range("E23").Select()
